$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 40; this shifts the existing rows 40-42 down to 41-43
$ws.Rows.Item(40).Insert()

# Populate the newly inserted row 40 with the new record (copy row 41's shared
# columns, which were previously row 40's, now shifted down)
$ws.Cells.Item(40, 1).Value = 11
$ws.Cells.Item(40, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(40, 3).Value = "Bíobío"
$ws.Cells.Item(40, 4).Value = 45077
$ws.Cells.Item(40, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(40, 5).Value = 8
$ws.Cells.Item(40, 6).Value = 100112022
$ws.Cells.Item(40, 7).Value = "Arveja Verde"
$ws.Cells.Item(40, 8).Value = "Perfection"
$ws.Cells.Item(40, 9).Value = "Primera"
$ws.Cells.Item(40, 10).Value = 100
$ws.Cells.Item(40, 11).Value = 30000
$ws.Cells.Item(40, 12).Value = 32000
$ws.Cells.Item(40, 13).Value = 31000
$ws.Cells.Item(40, 14).Value = "$/malla 25 kilos"
$ws.Cells.Item(40, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(40, 16).Value = 1240
$ws.Cells.Item(40, 17).Value = 25
$ws.Cells.Item(40, 18).Value = "Hortaliza"
